$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" column (C) date from 45207 to 45208 for rows 2-6
$ws.Range("C2:C6").Value = 45208
